# Inserts a new weekly price record for "Choclo" (Dulce o Americano,
# Región de Arica y Parinacota) on top of the existing records at row 70,
# shifting all subsequent rows (70-99) down by one (71-100).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 70; this shifts rows 70:99
# down to 71:100 and keeps every other row / formatting intact.
$ws.Rows.Item(70).Insert()

# Populate the newly inserted row 70 with the new record.
$ws.Range("A70").Value = 7
$ws.Range("B70").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C70").Value = "Ñuble"
$ws.Range("D70").Value = 44518
$ws.Range("E70").Value = 16
$ws.Range("F70").Value = 100112024
$ws.Range("G70").Value = "Choclo"
$ws.Range("H70").Value = "Dulce o Americano"
$ws.Range("I70").Value = "Primera"
$ws.Range("J70").Value = 60
$ws.Range("K70").Value = 16000
$ws.Range("L70").Value = 17000
$ws.Range("M70").Value = 16500
$ws.Range("N70").Value = '$/malla 60 unidades'
$ws.Range("O70").Value = "Región de Arica y Parinacota"
$ws.Range("P70").Value = 275
$ws.Range("Q70").Value = 60
$ws.Range("R70").Value = "Hortaliza"

# Match the date-number formatting used by the rest of column D.
$ws.Range("D70").NumberFormat = $ws.Range("D71").NumberFormat
